$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '="37.288.45"'
$ws.Range('D2').Copy()
$ws.Range('D2').PasteSpecial(-4163)
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Formula = '="2.060.50"'
$ws.Range('D3').Copy()
$ws.Range('D3').PasteSpecial(-4163)
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Formula = '="232.92"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('E6').Value = '  +2.50%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Formula = '="56.80"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  +0.89%  '
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('D10').Formula = '="57.97"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Formula = '="2.364.19"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Formula = '="14.60"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Formula = '="0.777"'
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('D17').Formula = '="5.15"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Formula = '="2.059.23"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('D19').Formula = '="37.208.39"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').Formula = '="6.34"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +7.44%  '
$ws.Range('D21').Formula = '="69.34"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('D22').Formula = '="0.0₃0809"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('D23').Formula = '="226.16"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  +1.43%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').Formula = '="2.44"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('E26').Value = '  -1.48%  '
$ws.Range('D27').Formula = '="166.60"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +6.72%  '
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('E30').Value = '  -1.36%  '
$ws.Range('D31').Formula = '="19.07"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('D33').Formula = '="4.44"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Formula = '="0.0618"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('D35').Formula = '="4.58"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  +5.76%  '
$ws.Range('D36').Formula = '="2.50"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('D40').Formula = '="5.68"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  -4.23%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').Formula = '="1.474.77"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').Formula = '="0.0938"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('D44').Formula = '="95.91"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +1.89%  '
$ws.Range('D45').Formula = '="0.0213"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +2.70%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').Formula = '="4.30"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  -1.31%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').Formula = '="1.17"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  +3.60%  '
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').Formula = '="15.08"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  -5.16%  '
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('E51').Value = '  +0.97%  '
